# Refresh the cryptocurrency market snapshot (price + 1h change) for every
# coin row on the sheet, and correct the ranking order of Dai /
# WrappedliquidstakedEther2.0 (rows 22-23 swap places).
#
# Every write goes through the same "store as literal text" pattern
# (NumberFormat "@" -> Value -> Style "Normal") so that values which look
# numeric (e.g. "29.895.34", "1.001", "0.000007787") round-trip exactly
# as text instead of being coerced into floating point numbers, while the
# cell keeps its original (default) style afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $value) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$updates = @(
    @("D2", "29.895.34"),
    @("E2", "  -0.31%  "),
    @("D3", "1.897.35"),
    @("E3", "  -0.08%  "),
    @("D4", "1.001"),
    @("E4", "  +0.02%  "),
    @("D5", "0.7925"),
    @("E5", "  -5.02%  "),
    @("D6", "244.19"),
    @("E6", "  +0.96%  "),
    @("E7", "  -0.04%  "),
    @("D8", "0.3168"),
    @("E8", "  -3.76%  "),
    @("D9", "25.43"),
    @("E9", "  -4.35%  "),
    @("D10", "0.07181"),
    @("E10", "  +1.83%  "),
    @("D11", "0.08115"),
    @("E11", "  +0.42%  "),
    @("D12", "5.638"),
    @("E12", "  +7.31%  "),
    @("D13", "0.7687"),
    @("E13", "  +0.89%  "),
    @("D14", "1.944.82"),
    @("E14", "  +2.58%  "),
    @("D15", "92.60"),
    @("E15", "  +0.33%  "),
    @("D16", "6.176"),
    @("E16", "  +5.11%  "),
    @("D17", "29.914.87"),
    @("E17", "  -0.24%  "),
    @("D18", "13.98"),
    @("E18", "  -0.99%  "),
    @("D19", "244.89"),
    @("E19", "  +0.18%  "),
    @("D20", "0.000007787"),
    @("E20", "  +0.34%  "),
    @("D21", "8.320"),
    @("E21", "  +19.29%  "),
    @("B22", "WrappedliquidstakedEther2.0"),
    @("C22", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"),
    @("D22", "2.157.45"),
    @("E22", "  +0.14%  "),
    @("B23", "Dai"),
    @("C23", "https://coinranking.com/coin/MoTuySvg7+dai-dai"),
    @("D23", "1.001"),
    @("E23", "  -0.02%  "),
    @("E24", "  +0.03%  "),
    @("D25", "0.1675"),
    @("E25", "  -3.74%  "),
    @("D26", "9.518"),
    @("E26", "  +2.71%  "),
    @("D27", "163.67"),
    @("E27", "  -1.38%  "),
    @("D28", "18.76"),
    @("E28", "  -0.81%  "),
    @("D29", "2.070"),
    @("E29", "  -1.25%  "),
    @("E30", "  +3.03%  "),
    @("E31", "  +2.42%  "),
    @("D32", "4.499"),
    @("E32", "  +4.90%  "),
    @("D33", "0.05635"),
    @("E33", "  -5.75%  "),
    @("D34", "4.101"),
    @("E34", "  +0.56%  "),
    @("D35", "1.281"),
    @("E35", "  +0.78%  "),
    @("D36", "0.7442"),
    @("E36", "  +1.66%  "),
    @("D37", "1.002"),
    @("E37", "  +0.24%  "),
    @("D38", "2.623"),
    @("E38", "  -3.74%  "),
    @("D39", "0.01934"),
    @("E39", "  +1.02%  "),
    @("D40", "2.787"),
    @("E40", "  +0.30%  "),
    @("D41", "1.167.29"),
    @("E41", "  +17.61%  "),
    @("D42", "74.97"),
    @("E42", "  +3.30%  "),
    @("D43", "0.4427"),
    @("E43", "  -0.46%  "),
    @("D44", "5.964"),
    @("E44", "  +1.81%  "),
    @("D45", "0.8543"),
    @("E45", "  +0.14%  "),
    @("D46", "104.74"),
    @("E46", "  +2.73%  "),
    @("E47", "  -0.08%  "),
    @("D48", "1.888"),
    @("E48", "  -0.45%  "),
    @("D49", "10.07"),
    @("E49", "  +2.56%  "),
    @("D50", "7.483"),
    @("E50", "  -1.00%  "),
    @("D51", "3.004"),
    @("E51", "  +10.08%  ")
)

foreach ($u in $updates) {
    Set-TextValue $u[0] $u[1]
}
